$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Insert two new blank columns before column C (old C->E, D->F, E->G, F->H, G->I).
#    Excel replicates the formatting of the column to the left (old B/C share style),
#    which is exactly what the target file expects for the two new columns.
$ws.Range("C1:D1").EntireColumn.Insert()

# 2. New header text in C1 ("Bank Uploaded"); D1 stays blank (already styled from insert).
$ws.Range("C1").Value = "Bank Uploaded"

# 3. Re-point the AutoFilter so it spans the new layout (B1:E18).
#    Toggling off first avoids the "no-op toggle" behaviour of Range.AutoFilter().
$ws.AutoFilterMode = $false
$ws.Range("B1:E18").AutoFilter()

# 4. Fix up the workbook-level hidden _FilterDatabase defined name so it matches
#    the new autofilter extent instead of the stale pre-insert range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $nm = $wb.Names.Item($i)
  if ($nm.Name -like "*_FilterDatabase*") {
    $nm.RefersTo = "=Sheet1!`$B`$1:`$E`$18"
  }
}

# 5. Rebuild the hyperlinks at their new (shifted) cell locations, in the same
#    order as before so relationship ids line up the same way.
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H3"), "http://akbl.com.pk/personal/consumer-products/master-card/discount-deal/")
$ws.Hyperlinks.Add($ws.Range("H4"), "https://www.bankalhabib.com/personal-banking/e-banking/discover.php")
$ws.Hyperlinks.Add($ws.Range("H5"), "http://www.bankalfalah.com/personal-banking/cards/dine-n-discount")
$ws.Hyperlinks.Add($ws.Range("H2"), "https://www.abl.com/personal-banking/credit-debit-cards/allied-visa-card-discounts-2017/")
$ws.Hyperlinks.Add($ws.Range("H8"), "http://www.hbl.com/hblcarddiscounts")
$ws.Hyperlinks.Add($ws.Range("H10"), "http://www.jsbl.com/products-services/card-products/credit-card/")
$ws.Hyperlinks.Add($ws.Range("H14"), "http://www.silkbank.com.pk/cc/")
$ws.Hyperlinks.Add($ws.Range("H16"), "https://www.sc.com/pk/credit-cards/the-good-life-privileges/")
$ws.Hyperlinks.Add($ws.Range("H17"), "http://summitbank.com.pk/index.php/electronic-banking/visa-debit-card/merchant-discounts/")
$ws.Hyperlinks.Add($ws.Range("H18"), "http://www.ubldirect.com/corporate/BankingServices/CardProducts/UBLDiscountAlliance.aspx")
$ws.Hyperlinks.Add($ws.Range("G15"), "mailto:service.quality@soneribank.com")
$ws.Hyperlinks.Add($ws.Range("G12"), "mailto:complaints.suggestion@nibpk.com")
$ws.Hyperlinks.Add($ws.Range("G13"), "mailto:samba.care@samba.com.pk")
$ws.Hyperlinks.Add($ws.Range("G11"), "mailto:info@mcb.com.pk")
$ws.Hyperlinks.Add($ws.Range("G7"), "mailto:info@habibmetro.com")
$ws.Hyperlinks.Add($ws.Range("G6"), "mailto:customercomplaint@faysalbank.com")
$ws.Hyperlinks.Add($ws.Range("H6"), "https://www.faysalbank.com/en/all-promotions/")
$ws.Hyperlinks.Add($ws.Range("H15"), "https://www.facebook.com/pg/SoneriBankPK/photos/?tab=album&album_id=560208060705453")

# 6. Column widths for the two newly inserted columns (closest representable
#    values - COM ColumnWidth quantises to 1/6 character, so we land on the
#    nearest reachable width to the authored 24.85546875 / 13.5703125).
$ws.Columns("C").ColumnWidth = 23.95
$ws.Columns("D").ColumnWidth = 12.6

# 7. Selection moves to F4 in the new layout.
$ws.Range("F4").Select()

Write-Host "done"
